$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the strain column (E) for rows 11-13 from "KN99alpha" to "TDY451"
$ws.Range("E11").Value = "TDY451"
$ws.Range("E12").Value = "TDY451"
$ws.Range("E13").Value = "TDY451"

# Match the new selection recorded in the workbook after the edit
$ws.Range("E12:E13").Select()
